# Updated cryptos list on Sat Sep  2 15:53:15 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns for each coin row, and
# fix the ordering of the WrappedEther / Polkadot rows (12 and 13).
#
# Note: several Price values look like plain numbers (e.g. "1.004",
# "215.90"). Assigning those bare strings to Range.Value would make Excel
# auto-convert them into real numbers (losing formatting / trailing
# zeros and changing the cell type). Prefixing with a leading apostrophe
# forces Excel to keep them as text, matching the original sheet, which
# stores every Price/Volume cell as a string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.983.28'
$ws.Range("E2").Value = '  -0.24%  '

$ws.Range("D3").Value = '1.645.67'
$ws.Range("E3").Value = '  -0.10%  '

$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  -0.42%  '

$ws.Range("D5").Value = '''215.90'
$ws.Range("E5").Value = '  +0.26%  '

$ws.Range("D6").Value = '''0.5066'
$ws.Range("E6").Value = '  -0.51%  '

$ws.Range("D7").Value = '''1.004'
$ws.Range("E7").Value = '  -0.41%  '

$ws.Range("D8").Value = '''0.2579'
$ws.Range("E8").Value = '  -0.21%  '

$ws.Range("D9").Value = '''0.06428'
$ws.Range("E9").Value = '  +0.11%  '

$ws.Range("D10").Value = '''19.70'
$ws.Range("E10").Value = '  +0.60%  '

$ws.Range("D11").Value = '''0.07767'
$ws.Range("E11").Value = '  +0.54%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.655.00'
$ws.Range("E12").Value = '  +0.41%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''4.277'
$ws.Range("E13").Value = '  +0.23%  '

$ws.Range("D14").Value = '1.871.29'
$ws.Range("E14").Value = '  -0.18%  '

$ws.Range("D15").Value = '''0.5469'
$ws.Range("E15").Value = '  +0.24%  '

$ws.Range("D16").Value = '0.0₅7958'
$ws.Range("E16").Value = '  -0.24%  '

$ws.Range("D17").Value = '''64.67'
$ws.Range("E17").Value = '  +1.44%  '

$ws.Range("D18").Value = '26.001.59'
$ws.Range("E18").Value = '  -0.29%  '

$ws.Range("D19").Value = '''1.005'
$ws.Range("E19").Value = '  -0.29%  '

$ws.Range("D20").Value = '''201.84'
$ws.Range("E20").Value = '  -2.25%  '

$ws.Range("D21").Value = '''4.399'
$ws.Range("E21").Value = '  +0.52%  '

$ws.Range("D22").Value = '''9.926'
$ws.Range("E22").Value = '  -0.90%  '

$ws.Range("D23").Value = '''6.001'
$ws.Range("E23").Value = '  +0.08%  '

$ws.Range("D24").Value = '''1.005'
$ws.Range("E24").Value = '  -0.35%  '

$ws.Range("D25").Value = '''1.879'
$ws.Range("E25").Value = '  +0.59%  '

$ws.Range("D26").Value = '''141.24'
$ws.Range("E26").Value = '  -1.34%  '

$ws.Range("D27").Value = '''0.1141'
$ws.Range("E27").Value = '  -2.05%  '

$ws.Range("D28").Value = '''6.838'
$ws.Range("E28").Value = '  -0.85%  '

$ws.Range("D29").Value = '''15.74'
$ws.Range("E29").Value = '  -0.28%  '

$ws.Range("D30").Value = '''1.246'
$ws.Range("E30").Value = '  +0.52%  '

$ws.Range("D31").Value = '''0.04942'
$ws.Range("E31").Value = '  -2.46%  '

$ws.Range("D32").Value = '''3.278'
$ws.Range("E32").Value = '  -1.42%  '

$ws.Range("D33").Value = '''3.220'
$ws.Range("E33").Value = '  -0.05%  '

$ws.Range("D34").Value = '''1.547'
$ws.Range("E34").Value = '  +0.12%  '

$ws.Range("D35").Value = '''2.375'
$ws.Range("E35").Value = '  +1.17%  '

$ws.Range("D36").Value = '''0.8958'
$ws.Range("E36").Value = '  -1.82%  '

$ws.Range("D37").Value = '''2.626'
$ws.Range("E37").Value = '  -0.77%  '

$ws.Range("D38").Value = '1.160.01'
$ws.Range("E38").Value = '  +1.23%  '

$ws.Range("D39").Value = '''0.5603'
$ws.Range("E39").Value = '  -1.54%  '

$ws.Range("D40").Value = '''0.01570'
$ws.Range("E40").Value = '  -0.08%  '

$ws.Range("D41").Value = '''1.006'
$ws.Range("E41").Value = '  -0.22%  '

$ws.Range("D42").Value = '''5.720'
$ws.Range("E42").Value = '  +1.12%  '

$ws.Range("D43").Value = '''0.8122'
$ws.Range("E43").Value = '  -1.36%  '

$ws.Range("D44").Value = '''99.81'
$ws.Range("E44").Value = '  +0.02%  '

$ws.Range("D45").Value = '1.781.17'
$ws.Range("E45").Value = '  -0.31%  '

$ws.Range("D46").Value = '0.0₈116'
$ws.Range("E46").Value = '  +2.10%  '

$ws.Range("E47").Value = '  -0.38%  '

$ws.Range("D48").Value = '''1.003'
$ws.Range("E48").Value = '  -0.66%  '

$ws.Range("D49").Value = '''54.96'
$ws.Range("E49").Value = '  -0.25%  '

$ws.Range("D50").Value = '''0.05054'
$ws.Range("E50").Value = '  -0.30%  '

$ws.Range("D51").Value = '''1.005'
$ws.Range("E51").Value = '  -0.41%  '
